# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the aggregated "全部类型" sheet, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (individual category)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 11892
$wsExhibit.Range("F4").Value = 19
$wsExhibit.Range("F5").Value = 220
$wsExhibit.Range("F8").Value = 11805
$wsExhibit.Range("F10").Value = 1173
$wsExhibit.Range("F11").Value = 100
$wsExhibit.Range("F12").Value = 59
$wsExhibit.Range("F13").Value = 1779
$wsExhibit.Range("F14").Value = 5852
$wsExhibit.Range("F15").Value = 125

# Sheet "全部类型" (all categories combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 11892
$wsAll.Range("F6").Value = 19
$wsAll.Range("F7").Value = 220
$wsAll.Range("F11").Value = 11805
$wsAll.Range("F13").Value = 1173
$wsAll.Range("F14").Value = 100
$wsAll.Range("F15").Value = 59
$wsAll.Range("F16").Value = 1779
$wsAll.Range("F18").Value = 5852
$wsAll.Range("F19").Value = 125
